$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.089.87"
$ws.Range("E2").Value = "  -2.41%  "
$ws.Range("D3").Value = "'3.520.99"
$ws.Range("E3").Value = "  -3.08%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'588.22"
$ws.Range("E5").Value = "  +0.73%  "
$ws.Range("D6").Value = "'170.18"
$ws.Range("E6").Value = "  -3.37%  "
$ws.Range("D7").Value = "'0.612"
$ws.Range("E7").Value = "  -0.36%  "
$ws.Range("D8").Value = "'3.513.66"
$ws.Range("E8").Value = "  -3.08%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").Value = "'0.189"
$ws.Range("E10").Value = "  -4.22%  "
$ws.Range("D11").Value = "'6.82"
$ws.Range("E11").Value = "  -0.27%  "
$ws.Range("E12").Value = "  -5.05%  "
$ws.Range("D13").Value = "'47.42"
$ws.Range("E13").Value = "  -2.40%  "
$ws.Range("E14").Value = "  -3.18%  "
$ws.Range("D15").Value = "'4.088.82"
$ws.Range("E15").Value = "  -3.16%  "
$ws.Range("E16").Value = "  -5.82%  "
$ws.Range("D17").Value = "'614.61"
$ws.Range("E17").Value = "  -8.86%  "
$ws.Range("D18").Value = "'3.531.48"
$ws.Range("E18").Value = "  -2.91%  "
$ws.Range("D19").Value = "'69.156.22"
$ws.Range("E19").Value = "  -2.37%  "
$ws.Range("E20").Value = "  -1.33%  "
$ws.Range("D21").Value = "'17.41"
$ws.Range("E21").Value = "  -2.19%  "
$ws.Range("E22").Value = "  -3.54%  "
$ws.Range("D23").Value = "'0.886"
$ws.Range("E23").Value = "  -5.99%  "
$ws.Range("D24").Value = "'15.80"
$ws.Range("E24").Value = "  -8.00%  "
$ws.Range("D25").Value = "'96.59"
$ws.Range("E25").Value = "  -3.47%  "
$ws.Range("D26").Value = "'3.84"
$ws.Range("E26").Value = "  -2.31%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").Value = "'2.62"
$ws.Range("E28").Value = "  -6.37%  "
$ws.Range("D29").Value = "'9.22"
$ws.Range("E29").Value = "  -6.04%  "
$ws.Range("D30").Value = "'32.62"
$ws.Range("E30").Value = "  -5.88%  "
$ws.Range("D31").Value = "'8.53"
$ws.Range("E31").Value = "  -6.78%  "
$ws.Range("D32").Value = "'3.13"
$ws.Range("E32").Value = "  -4.94%  "
$ws.Range("D33").Value = "'1.32"
$ws.Range("E33").Value = "  -5.13%  "
$ws.Range("D34").Value = "'6.93"
$ws.Range("E34").Value = "  -8.54%  "
$ws.Range("D35").Value = "'615.75"
$ws.Range("E35").Value = "  +6.94%  "
$ws.Range("D36").Value = "'10.75"
$ws.Range("E36").Value = "  -3.18%  "
$ws.Range("D37").Value = "'3.48"
$ws.Range("E37").Value = "  -12.63%  "
$ws.Range("E38").Value = "  -4.64%  "
$ws.Range("D39").Value = "'57.13"
$ws.Range("E39").Value = "  -2.29%  "
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("D41").Value = "'0.0444"
$ws.Range("E41").Value = "  -1.71%  "
$ws.Range("E42").Value = "  -3.24%  "
$ws.Range("D43").Value = "'3.401.51"
$ws.Range("E43").Value = "  -4.39%  "
$ws.Range("D44").Value = "'0.326"
$ws.Range("E44").Value = "  -5.55%  "
$ws.Range("D45").Value = "'32.84"
$ws.Range("E45").Value = "  -4.29%  "
$ws.Range("D46").Value = "'0.0₃0699"
$ws.Range("E46").Value = "  -4.66%  "
$ws.Range("D47").Value = "'2.53"
$ws.Range("E47").Value = "  -5.63%  "
$ws.Range("E48").Value = "  -6.20%  "
$ws.Range("E49").Value = "  -3.21%  "
$ws.Range("D50").Value = "'134.06"
$ws.Range("E50").Value = "  -2.69%  "
$ws.Range("D51").Value = "'5.53"
$ws.Range("E51").Value = "  +10.30%  "
